# Update database: drop oldest quarter (1399/06), shift remaining quarters
# left by one column, and append the newest quarter (1401/12) in column M.
# Also refresh the "read_price" derived figures for the existing quarters
# (values shift from the next column) and populate the new quarter's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (period headers): shift text left one column, set the new last quarter ---
$ws.Range("D8").Value = $ws.Range("E8").Value2
$ws.Range("E8").Value = $ws.Range("F8").Value2
$ws.Range("F8").Value = $ws.Range("G8").Value2
$ws.Range("G8").Value = $ws.Range("H8").Value2
$ws.Range("H8").Value = $ws.Range("I8").Value2
$ws.Range("I8").Value = $ws.Range("J8").Value2
$ws.Range("J8").Value = $ws.Range("K8").Value2
$ws.Range("K8").Value = $ws.Range("L8").Value2
$ws.Range("L8").Value = $ws.Range("M8").Value2
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9 (publish dates): shift text left one column, set the new last date ---
# Dates are stored as plain text (not real dates), so use a leading apostrophe
# to stop Excel from auto-converting them into date serial numbers.
$ws.Range("D9").Value = "'" + $ws.Range("E9").Text
$ws.Range("E9").Value = "'" + $ws.Range("F9").Text
$ws.Range("F9").Value = "'" + $ws.Range("G9").Text
$ws.Range("G9").Value = "'" + $ws.Range("H9").Text
$ws.Range("H9").Value = "'" + $ws.Range("I9").Text
$ws.Range("I9").Value = "'1402-02-10 (7)"
$ws.Range("J9").Value = "'" + $ws.Range("K9").Text
$ws.Range("K9").Value = "'" + $ws.Range("L9").Text
$ws.Range("L9").Value = "'" + $ws.Range("M9").Text
$ws.Range("M9").Value = "'1402-02-10"

# --- Numeric data rows: shift left by one column, then fill column M with the new quarter's figures ---
$newValues = @{
    11 = 29289
    12 = -12272
    13 = 17017
    14 = -4685
    16 = 2389
    17 = 14720
    19 = 2607
    20 = 17327
    21 = 924
    22 = 18251
    24 = 18251
    26 = 3938
}

foreach ($row in @(11,12,13,14,16,17,19,20,21,22,24,26)) {
    $ws.Range("D" + $row).Value = $ws.Range("E" + $row).Value2
    $ws.Range("E" + $row).Value = $ws.Range("F" + $row).Value2
    $ws.Range("F" + $row).Value = $ws.Range("G" + $row).Value2
    $ws.Range("G" + $row).Value = $ws.Range("H" + $row).Value2
    $ws.Range("H" + $row).Value = $ws.Range("I" + $row).Value2
    $ws.Range("I" + $row).Value = $ws.Range("J" + $row).Value2
    $ws.Range("J" + $row).Value = $ws.Range("K" + $row).Value2
    $ws.Range("K" + $row).Value = $ws.Range("L" + $row).Value2
    $ws.Range("L" + $row).Value = $ws.Range("M" + $row).Value2
    $ws.Range("M" + $row).Value = $newValues[$row]
}
